$wb = $excel.ActiveWorkbook

# Sheets involved in the localization-status report.
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# The handoff status moved from "Ready for handoff" to "In Translation"
# for both tracked files, on every sheet that surfaces the Status column.
$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) / F (de-de) hold the per-locale status.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Per-locale detail sheets: column C is "Status".
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# The shorter status text means the Status columns no longer need to be as
# wide, so re-fit them to the new content width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
